$d = $word.ActiveDocument
$p = $d.Paragraphs.Last

# Paragraph 1: ilvl=0
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Done in R"
$p.Range.ListFormat.ListLevelNumber = 1

# Paragraph 2: ilvl=0
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Book series analysis"
$p.Range.ListFormat.ListLevelNumber = 1

# Paragraph 3: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Both series express a very erratic behavior, although they seem to be trending up in both cases. Since the frequency of the data is 1, it is hard to discuss any type of seasonality. A cycle also seems not be clear from the data."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 4: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "For the paperback series, we expect the next four days to sell 207 books. For the hardcover series, we expect to sell 240 books in the next four days."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 5: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Paperback RMSE: 33.6 / Hardcover RMSE: 31.9"
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 6: ilvl=0
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Using holt’s linear method."
$p.Range.ListFormat.ListLevelNumber = 1

# Paragraph 7: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Done in R."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 8: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "RMSE for paperback using Holt’s method: 31.1/ RMSE for hardcover using Holt’s method: 27.2. The RMSE is smaller using Holt’s method, though this is expected since it is a more complex model. The Holt’s method seems to be more adequate, because both series appear to be trending up. However, the estimates of beta for both series are really small, meaning that there is little evidence that the trend is changing over time."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 9: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "The Holt’s forecasts seem to be better for Hardcover series. For the paperback series"
$p.Range.ListFormat.ListLevelNumber = 2
$rEnd = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$rEnd.InsertAfter(", the reduction in RMSE was not as pronounced, so there isn’t much evidence that the more complicated model provides a better fit.")

# Paragraph 10: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "For both SES and Holt’s linear method, the 95% CIs produced by using the RMSE are narrower than the ones produced by the function"
$p.Range.ListFormat.ListLevelNumber = 2
$rEnd = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$rEnd.InsertAfter("s")
$rEnd = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$rEnd.InsertAfter(".")

# Paragraph 11: ilvl=0
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Eggs dataset."
$p.Range.ListFormat.ListLevelNumber = 1

# Paragraph 12: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "As expected, the forecasts without dumping have a linear downward trend. The estimate for beta is small, meaning that the variation in the trend is almost negligible as time goes by. The estimate for alpha is very large, suggesting that the method relies considerably on new information introduced."
$p.Range.ListFormat.ListLevelNumber = 2
$rEnd = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$rEnd.InsertAfter(" Using a damping parameter of 0.85 makes the forecasts barely change. Introducing a box-cox transformation")
$rEnd = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$rEnd.InsertAfter(" does not change the overall trends.")

# Paragraph 13: ilvl=0
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Retail dataset."
$p.Range.ListFormat.ListLevelNumber = 1

# Paragraph 14: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Because the seasonal variability is not constant."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 15: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Using a damp parameter prevents the forecasts from growing indefinitely into the future."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 16: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Using time-series cross validation, the RMSE of the method without dumping was smaller."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 17: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "They do not look like White-Noise. The Ljung-Box text rejects the null hypothesis of no residual autocorrelation, while the correlogram shows strong correlations especially at lags multiples of 6."
$p.Range.ListFormat.ListLevelNumber = 2

# Paragraph 18: ilvl=1
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "The model with a damp term performs better than the seasonal naïve method."
$p.Range.ListFormat.ListLevelNumber = 2
